$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.10402464666999
$ws.Range("C2").Value = 6.605890326123556
$ws.Range("D2").Value = 10.06068743220359
$ws.Range("F2").Value = 35.84741385976312
$ws.Range("G2").Value = 38.38553369384999
$ws.Range("H2").Value = 16.661634794818
$ws.Range("I2").Value = 23.78251054217777
$ws.Range("J2").Value = 11.40713408812248
$ws.Range("K2").Value = 11.5094742519409
$ws.Range("N2").Value = 19.47155882203474
$ws.Range("B3").Value = 10.84240845027909
$ws.Range("C3").Value = 6.408871746482109
$ws.Range("D3").Value = 9.976600997507543
$ws.Range("F3").Value = 35.846223594061
$ws.Range("G3").Value = 38.38502775382708
$ws.Range("H3").Value = 16.70479343571246
$ws.Range("I3").Value = 23.85757545181176
$ws.Range("J3").Value = 11.39008536600348
$ws.Range("K3").Value = 11.33472328269805
$ws.Range("N3").Value = 19.53581380096228
$ws.Range("B4").Value = 10.68097425995302
$ws.Range("C4").Value = 6.286389066424225
$ws.Range("D4").Value = 9.92671928208437
$ws.Range("F4").Value = 35.85439256951782
$ws.Range("G4").Value = 38.39606701838809
$ws.Range("H4").Value = 16.73413182004761
$ws.Range("I4").Value = 23.90823762047877
$ws.Range("J4").Value = 11.38192175820146
$ws.Range("K4").Value = 11.22851366106375
$ws.Range("N4").Value = 19.57706484800984
$ws.Range("B5").Value = 10.61508913256396
$ws.Range("C5").Value = 6.236179453921483
$ws.Range("D5").Value = 9.90685083142362
$ws.Range("F5").Value = 35.85995850767058
$ws.Range("G5").Value = 38.40341509033405
$ws.Range("H5").Value = 16.7468004504013
$ws.Range("I5").Value = 23.9300303534531
$ws.Range("J5").Value = 11.37917673372885
$ws.Range("K5").Value = 11.18556102633872
$ws.Range("N5").Value = 19.59432853567608
$ws.Range("B6").Value = 10.60414600426195
$ws.Range("C6").Value = 6.2278267893505
$ws.Range("D6").Value = 9.903579944600347
$ws.Range("F6").Value = 35.86101778394463
$ws.Range("G6").Value = 38.4048071598221
$ws.Range("H6").Value = 16.74894710496108
$ws.Range("I6").Value = 23.93371826092495
$ws.Range("J6").Value = 11.37875611584672
$ws.Range("K6").Value = 11.17845022668236
$ws.Range("N6").Value = 19.59722259155759
$ws.Range("B7").Value = 10.68008597512994
$ws.Range("C7").Value = 6.28571301032061
$ws.Range("D7").Value = 9.926449446894498
$ws.Range("F7").Value = 35.8544585784057
$ws.Range("G7").Value = 38.39615458724838
$ws.Range("H7").Value = 16.73429978787878
$ws.Range("I7").Value = 23.90852688171285
$ws.Range("J7").Value = 11.38188237985553
$ws.Range("K7").Value = 11.22793298362089
$ws.Range("N7").Value = 19.57729583392528
$ws.Range("B8").Value = 11.0140456112162
$ws.Range("C8").Value = 6.538320654518006
$ws.Range("D8").Value = 10.03134351443181
$ws.Range("F8").Value = 35.84515615298688
$ws.Range("G8").Value = 38.38300145948556
$ws.Range("H8").Value = 16.67592620707044
$ws.Range("I8").Value = 23.8074427054968
$ws.Range("J8").Value = 11.40077889766457
$ws.Range("K8").Value = 11.44902731173841
$ws.Range("N8").Value = 19.49334154598113
$ws.Range("B9").Value = 11.65830805114891
$ws.Range("C9").Value = 7.018290966166156
$ws.Range("D9").Value = 10.25002298250521
$ws.Range("F9").Value = 35.89751470194918
$ws.Range("G9").Value = 38.4473757914193
$ws.Range("H9").Value = 16.58401703854962
$ws.Range("I9").Value = 23.64558993566007
$ws.Range("J9").Value = 11.4560019240184
$ws.Range("K9").Value = 11.88881149526043
$ws.Range("N9").Value = 19.34291176025022
$ws.Range("B10").Value = 12.1197719052258
$ws.Range("C10").Value = 7.357402628283151
$ws.Range("D10").Value = 10.41737150087737
$ws.Range("F10").Value = 35.97891947421144
$ws.Range("G10").Value = 38.54965708354825
$ws.Range("H10").Value = 16.53029391461566
$ws.Range("I10").Value = 23.54898003209166
$ws.Range("J10").Value = 11.5074662898642
$ws.Range("K10").Value = 12.21251195138464
$ws.Range("N10").Value = 19.24095920079437
$ws.Range("B11").Value = 12.3260905011904
$ws.Range("C11").Value = 7.507979429078353
$ws.Range("D11").Value = 10.49468269676926
$ws.Range("F11").Value = 36.02522375104964
$ws.Range("G11").Value = 38.60808375214962
$ws.Range("H11").Value = 16.50886053021581
$ws.Range("I11").Value = 23.50989934327642
$ws.Range("J11").Value = 11.53319581357081
$ws.Range("K11").Value = 12.35924156741034
$ws.Range("N11").Value = 19.1964193627953
$ws.Range("B12").Value = 12.40362142933437
$ws.Range("C12").Value = 7.564413532170583
$ws.Range("D12").Value = 10.5241054100361
$ws.Range("F12").Value = 36.04408439491807
$ws.Range("G12").Value = 38.63191182426542
$ws.Range("H12").Value = 16.50117716557246
$ws.Range("I12").Value = 23.49580258746496
$ws.Range("J12").Value = 11.54326746837245
$ws.Range("K12").Value = 12.41467610371495
$ws.Range("N12").Value = 19.17981627905209
$ws.Range("B13").Value = 12.38695152669295
$ws.Range("C13").Value = 7.552286334345632
$ws.Range("D13").Value = 10.51776255622434
$ws.Range("F13").Value = 36.0399635646179
$ws.Range("G13").Value = 38.62670442326426
$ws.Range("H13").Value = 16.50281264601156
$ws.Range("I13").Value = 23.49880730745793
$ws.Range("J13").Value = 11.54108383222344
$ws.Range("K13").Value = 12.40274384372494
$ws.Range("N13").Value = 19.1833803646183
$ws.Range("B14").Value = 12.33248142373231
$ws.Range("C14").Value = 7.512634342596153
$ws.Range("D14").Value = 10.49710051738193
$ws.Range("F14").Value = 36.02674889480198
$ws.Range("G14").Value = 38.61001002115145
$ws.Range("H14").Value = 16.50821973298186
$ws.Range("I14").Value = 23.5087255076357
$ws.Range("J14").Value = 11.5340178617553
$ws.Range("K14").Value = 12.36380508079939
$ws.Range("N14").Value = 19.1950481495484
$ws.Range("B15").Value = 12.29903685761615
$ws.Range("C15").Value = 7.488268466549666
$ws.Range("D15").Value = 10.48446281853427
$ws.Range("F15").Value = 36.0188270012415
$ws.Range("G15").Value = 38.60000574789854
$ws.Range("H15").Value = 16.51158814110105
$ws.Range("I15").Value = 23.51489221243421
$ws.Range("J15").Value = 11.52973237112106
$ws.Range("K15").Value = 12.33993563699492
$ws.Range("N15").Value = 19.2022292425917
$ws.Range("B16").Value = 12.10620907490791
$ws.Range("C16").Value = 7.347482984565605
$ws.Range("D16").Value = 10.41234083752782
$ws.Range("F16").Value = 35.97607939282796
$ws.Range("G16").Value = 38.54607757484221
$ws.Range("H16").Value = 16.53175518548502
$ws.Range("I16").Value = 23.55163219400318
$ws.Range("J16").Value = 11.5058310422688
$ws.Range("K16").Value = 12.20290767564369
$ws.Range("N16").Value = 19.24390688056742
$ws.Range("B17").Value = 11.986934222364
$ws.Range("C17").Value = 7.260129797707723
$ws.Range("D17").Value = 10.36838275251902
$ws.Range("F17").Value = 35.95222547196568
$ws.Range("G17").Value = 38.51603703446516
$ws.Range("H17").Value = 16.54489739134096
$ws.Range("I17").Value = 23.57541935633517
$ws.Range("J17").Value = 11.49175872612532
$ws.Range("K17").Value = 12.11867273788582
$ws.Range("N17").Value = 19.26994485735582
$ws.Range("B18").Value = 11.91799565310376
$ws.Range("C18").Value = 7.20954251225335
$ws.Range("D18").Value = 10.34321204347087
$ws.Range("F18").Value = 35.93937865844912
$ws.Range("G18").Value = 38.49987905985508
$ws.Range("H18").Value = 16.55273929125008
$ws.Range("I18").Value = 23.58955913568658
$ws.Range("J18").Value = 11.48388321556867
$ws.Range("K18").Value = 12.07017685991233
$ws.Range("N18").Value = 19.28509437183984
$ws.Range("B19").Value = 11.89459924338803
$ws.Range("C19").Value = 7.192357233090791
$ws.Range("D19").Value = 10.33470976607483
$ws.Range("F19").Value = 35.93517913827127
$ws.Range("G19").Value = 38.49460090266673
$ws.Range("H19").Value = 16.55544297838869
$ws.Range("I19").Value = 23.59442521858637
$ws.Range("J19").Value = 11.48125437532371
$ws.Range("K19").Value = 12.05375072117404
$ws.Range("N19").Value = 19.29025351546336
$ws.Range("B20").Value = 11.99966648055333
$ws.Range("C20").Value = 7.269464723556592
$ws.Range("D20").Value = 10.37305065961076
$ws.Range("F20").Value = 35.95467442342221
$ws.Range("G20").Value = 38.51911897771457
$ws.Range("H20").Value = 16.54346910109908
$ws.Range("I20").Value = 23.57283975075884
$ws.Range("J20").Value = 11.49323416704279
$ws.Range("K20").Value = 12.12764486301231
$ws.Range("N20").Value = 19.2671551580049
$ws.Range("B21").Value = 12.34849742984257
$ws.Range("C21").Value = 7.524297417109025
$ws.Range("D21").Value = 10.503165667347
$ws.Range("F21").Value = 36.03059443535048
$ws.Range("G21").Value = 38.61486742226513
$ws.Range("H21").Value = 16.50661978307249
$ws.Range("I21").Value = 23.50579321498287
$ws.Range("J21").Value = 11.53608443586168
$ws.Range("K21").Value = 12.3752462514846
$ws.Range("N21").Value = 19.19161390534553
$ws.Range("B22").Value = 12.57296369483239
$ws.Range("C22").Value = 7.68740725474076
$ws.Range("D22").Value = 10.58904712483834
$ws.Range("F22").Value = 36.0879383666405
$ws.Range("G22").Value = 38.68736662493448
$ws.Range("H22").Value = 16.4850607672858
$ws.Range("I22").Value = 23.46606870506753
$ws.Range("J22").Value = 11.56600141257911
$ws.Range("K22").Value = 12.53629400652379
$ws.Range("N22").Value = 19.14377681103671
$ws.Range("B23").Value = 12.45350808015608
$ws.Range("C23").Value = 7.600684081081782
$ws.Range("D23").Value = 10.54314116530989
$ws.Range("F23").Value = 36.056628688542
$ws.Range("G23").Value = 38.64776774969442
$ws.Range("H23").Value = 16.49633601105703
$ws.Range("I23").Value = 23.48689505927204
$ws.Range("J23").Value = 11.54986095135268
$ws.Range("K23").Value = 12.45042778654703
$ws.Range("N23").Value = 19.16916846328941
$ws.Range("B24").Value = 11.99391136011482
$ws.Range("C24").Value = 7.265245541691965
$ws.Range("D24").Value = 10.37093998124171
$ws.Range("F24").Value = 35.95356455084398
$ws.Range("G24").Value = 38.51772216443907
$ws.Range("H24").Value = 16.54411393942513
$ws.Range("I24").Value = 23.57400454453745
$ws.Range("J24").Value = 11.49256645060157
$ws.Range("K24").Value = 12.1235887732084
$ws.Range("N24").Value = 19.26841582046516
$ws.Range("B25").Value = 11.48573339505858
$ws.Range("C25").Value = 6.890544275511578
$ws.Range("D25").Value = 10.1895993685566
$ws.Range("F25").Value = 35.87579684522795
$ws.Range("G25").Value = 38.42029706366613
$ws.Range("H25").Value = 16.60646035996067
$ws.Range("I25").Value = 23.68546783406977
$ws.Range("J25").Value = 11.43913486467098
$ws.Range("K25").Value = 11.76950625880577
$ws.Range("N25").Value = 19.3820957921879
